# Append a new data row (row 99) to each of the 4 worksheets,
# mirroring the structure/format of the existing rows (A: datetime
# with the "YYYY-MM-DD HH:MM:SS" number format, B-E: text, F-I: numbers).

$wb = $excel.ActiveWorkbook

$newRowDate = [double]"45885.43790509259"

# --- Sheet 1: DE_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$r = 99
$ws.Cells.Item($r, 1).Value = $newRowDate
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x7c"
$ws.Cells.Item($r, 3).Value = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
$ws.Cells.Item($r, 4).Value = "0x01,0x24"
$ws.Cells.Item($r, 5).Value = "0x14"
$ws.Cells.Item($r, 6).Value = 380
$ws.Cells.Item($r, 7).Value = [double]"7.598631275147109e+23"
$ws.Cells.Item($r, 8).Value = 292
$ws.Cells.Item($r, 9).Value = 14

# --- Sheet 2: DE_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$r = 99
$ws.Cells.Item($r, 1).Value = $newRowDate
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x01,0x7c"
$ws.Cells.Item($r, 3).Value = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
$ws.Cells.Item($r, 4).Value = "0x01,0x28"
$ws.Cells.Item($r, 5).Value = "0xe"
$ws.Cells.Item($r, 6).Value = 380
$ws.Cells.Item($r, 7).Value = [double]"5.68432987514711e+23"
$ws.Cells.Item($r, 8).Value = 296
$ws.Cells.Item($r, 9).Value = 14

# --- Sheet 3: DE_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$r = 99
$ws.Cells.Item($r, 1).Value = $newRowDate
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x00,0x82"
$ws.Cells.Item($r, 3).Value = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x00,0x74"
$ws.Cells.Item($r, 5).Value = "0x7"
$ws.Cells.Item($r, 6).Value = 130
$ws.Cells.Item($r, 7).Value = [double]"5.68631262647114e+23"
$ws.Cells.Item($r, 8).Value = 116
$ws.Cells.Item($r, 9).Value = 7

# --- Sheet 4: DE_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$r = 99
$ws.Cells.Item($r, 1).Value = $newRowDate
$ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item($r, 2).Value = "0x00,0x82"
$ws.Cells.Item($r, 3).Value = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
$ws.Cells.Item($r, 4).Value = "0x00,0x73"
$ws.Cells.Item($r, 5).Value = "0x3"
$ws.Cells.Item($r, 6).Value = 130
$ws.Cells.Item($r, 7).Value = [double]"9.85046333984776e+23"
$ws.Cells.Item($r, 8).Value = 115
$ws.Cells.Item($r, 9).Value = 3

Write-Output "Appended row 99 to all 4 worksheets"
